# Insert a new weekly record as row 71 on the "Camote" sheet.
# This pushes the existing rows 71-89 down to 72-90 (dimension grows from
# A1:R89 to A1:R90), and populates the new row 71 with the latest entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 71:89 down by inserting a fresh row at 71.
$ws.Rows.Item(71).Insert()

# Populate the newly inserted row 71 with the new data record.
$ws.Cells.Item(71, 1).Value  = 10
$ws.Cells.Item(71, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(71, 3).Value  = "La Araucanía"
$ws.Cells.Item(71, 4).Value  = 44722
$ws.Cells.Item(71, 5).Value  = 9
$ws.Cells.Item(71, 6).Value  = 100114002
$ws.Cells.Item(71, 7).Value  = "Camote"
$ws.Cells.Item(71, 8).Value  = "Sin especificar"
$ws.Cells.Item(71, 9).Value  = "Primera"
$ws.Cells.Item(71, 10).Value = 20
$ws.Cells.Item(71, 11).Value = 20000
$ws.Cells.Item(71, 12).Value = 20000
$ws.Cells.Item(71, 13).Value = 20000
$ws.Cells.Item(71, 14).Value = "`$/malla 20 kilos"
$ws.Cells.Item(71, 15).Value = "Perú"
$ws.Cells.Item(71, 16).Value = 1000
$ws.Cells.Item(71, 17).Value = 20
$ws.Cells.Item(71, 18).Value = "Hortaliza"
